# Updates stack-trace line numbers / class names embedded in the
# bold red error text (JUnit/M2Doc moving from 3.1.0 to 3.1.1).
$d = $word.ActiveDocument
$replacements = 0

$found0 = $d.Content.Find.Execute("`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:192)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)", 2)
if ($found0) { $replacements = $replacements + 1 }
$found1 = $d.Content.Find.Execute("`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109)", 2)
if ($found1) { $replacements = $replacements + 1 }
$found2 = $d.Content.Find.Execute("`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:586)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:587)", 2)
if ($found2) { $replacements = $replacements + 1 }
$found3 = $d.Content.Find.Execute("`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1464)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1467)", 2)
if ($found3) { $replacements = $replacements + 1 }
$found4 = $d.Content.Find.Execute("`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:296)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:297)", 2)
if ($found4) { $replacements = $replacements + 1 }
$found5 = $d.Content.Find.Execute("`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:281)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:805)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:420)`n`tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:282)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:845)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:514)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:421)`n`tat sun.reflect.GeneratedMethodAccessor73.invoke(Unknown Source)", 2)
if ($found5) { $replacements = $replacements + 1 }
$found6 = $d.Content.Find.Execute("`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:50)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:59)", 2)
if ($found6) { $replacements = $replacements + 1 }
$found7 = $d.Content.Find.Execute("`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)", 2)
if ($found7) { $replacements = $replacements + 1 }
$found8 = $d.Content.Find.Execute("`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.BlockJUnit4ClassRunner`$1.evaluate(BlockJUnit4ClassRunner.java:100)`n`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)", 2)
if ($found8) { $replacements = $replacements + 1 }
$found9 = $d.Content.Find.Execute("`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)", 2)
if ($found9) { $replacements = $replacements + 1 }
$found10 = $d.Content.Find.Execute("`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)", 2)
if ($found10) { $replacements = $replacements + 1 }
$found11 = $d.Content.Find.Execute("`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)", 2)
if ($found11) { $replacements = $replacements + 1 }
$found12 = $d.Content.Find.Execute("`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)", $true, $false, $false, $false, $false, $true, 1, $false, "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)", 2)
if ($found12) { $replacements = $replacements + 1 }

Write-Host "Applied $replacements stack-trace text replacements"
